# Mother-In-Law House Expenses - add missing Friday 24/10/2025 labor entries
# and roll the corrected totals through every dependent summary sheet.
#
# Notes on technique:
#  - This workbook stores every number/label as a literal value (no live
#    formulas anywhere), so each derived total has to be poked by hand on
#    every sheet that shows it, exactly like the original diff does.
#  - A handful of the text cells look like a percentage ("151.63%"). Plain
#    `.Value = "151.63%"` gets auto-detected as a real percentage number
#    (0.1516 formatted as 0.00%), which would change the cell's type/style.
#    Writing it through a Text-formatted cell and then resetting the style
#    back to Normal keeps it as plain literal text with the default style,
#    matching the source file.
#  - New rows are cloned from an existing row with the same PAID/UNPAID
#    look via Copy + PasteSpecial(xlPasteFormats) so they pick up the exact
#    same style indices the workbook already defines, instead of Excel
#    synthesizing brand-new styles.

$wb = $excel.ActiveWorkbook

function Set-LiteralText {
    param($range, [string]$text)
    # Force a value to be stored as plain text even if it "looks like" a
    # number/percentage/date, then drop the temporary Text number format so
    # the cell is left with the default (General) style, same as the source.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# Sheet: Home Summary
# ---------------------------------------------------------------------
$home = $wb.Worksheets.Item("Home Summary")

$home.Range("B5").Value = "KES 1,516,311"
$home.Range("B6").Value = "KES -516,311"
Set-LiteralText $home.Range("B7") "151.63%"
$home.Range("B8").Value = "KES 17,399"

$home.Range("B12").Value = "KES 38,900"
$home.Range("B13").Value = "KES 120,200"
$home.Range("B14").Value = "KES 1,636,511"
Set-LiteralText $home.Range("B15") "163.65%"
$home.Range("B16").Value = "KES -636,511"

$home.Range("B19").Value = "KES 170,405"
$home.Range("B20").Value = "KES 1,806,916"
$home.Range("B21").Value = "KES 806,916"

$home.Range("B25").Value = 504050
$home.Range("C25").Value = 6320.5
$home.Range("D25").Value = 510370.5
Set-LiteralText $home.Range("E25") "51.04%"

# ---------------------------------------------------------------------
# Sheet: Daily Expenses - append the missing Fri 24/10 + Sat 25/10 rows
# ---------------------------------------------------------------------
$daily = $wb.Worksheets.Item("Daily Expenses")

# Clone the look of an existing "UNPAID" labor block (rows 455-458, the
# Sun 19/10 entries) onto the ten new rows, then overwrite the text/values.
$daily.Range("A455:I458").Copy()
$daily.Range("A489:I492").PasteSpecial(-4122)
$daily.Range("A455:I458").Copy()
$daily.Range("A493:I496").PasteSpecial(-4122)
$daily.Range("A479:I479").Copy()
$daily.Range("A497:I497").PasteSpecial(-4122)
$daily.Range("A488:I488").Copy()
$daily.Range("A498:I498").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$dailyRows = @(
    @(489, "24/10/2025", "Labor Costs", "Daily Labor", "Jack - UNPAID", 1500, 0, 0, "Worker", "UNPAID"),
    @(490, "24/10/2025", "Labor Costs", "Daily Labor", "Fundi 1 - UNPAID", 1300, 0, 0, "Worker", "UNPAID"),
    @(491, "24/10/2025", "Labor Costs", "Daily Labor", "Fundi 2 - UNPAID", 1300, 0, 0, "Worker", "UNPAID"),
    @(492, "24/10/2025", "Labor Costs", "Daily Labor", "2 helpers @ 600 each - UNPAID", 1200, 0, 0, "Worker", "UNPAID"),
    @(493, "25/10/2025", "Labor Costs", "Daily Labor", "Jack - UNPAID", 1500, 0, 0, "Worker", "UNPAID"),
    @(494, "25/10/2025", "Labor Costs", "Daily Labor", "Fundi 1 - UNPAID", 1300, 0, 0, "Worker", "UNPAID"),
    @(495, "25/10/2025", "Labor Costs", "Daily Labor", "Fundi 2 - UNPAID", 1300, 0, 0, "Worker", "UNPAID"),
    @(496, "25/10/2025", "Labor Costs", "Daily Labor", "2 helpers @ 600 each - UNPAID", 1200, 0, 0, "Worker", "UNPAID"),
    @(497, "25/10/2025", "Transport & Logistics", "Worker Transport", "Transport - UNPAID", 600, 0, 0, "Local Transport", "UNPAID"),
    @(498, "25/10/2025", "Building Materials", "Paint & Finishes", "5kg red oxide @ 200", 1000, 10, 1010, "Hardware Store", "PAID")
)
foreach ($entry in $dailyRows) {
    $r = $entry[0]
    $daily.Cells.Item($r, 1).Value = $entry[1]
    $daily.Cells.Item($r, 2).Value = $entry[2]
    $daily.Cells.Item($r, 3).Value = $entry[3]
    $daily.Cells.Item($r, 4).Value = $entry[4]
    $daily.Cells.Item($r, 5).Value = $entry[5]
    $daily.Cells.Item($r, 6).Value = $entry[6]
    $daily.Cells.Item($r, 7).Value = $entry[7]
    $daily.Cells.Item($r, 8).Value = $entry[8]
    $daily.Cells.Item($r, 9).Value = $entry[9]
}

# ---------------------------------------------------------------------
# Sheet: M-Pesa Fees
# ---------------------------------------------------------------------
$mpesa = $wb.Worksheets.Item("M-Pesa Fees")
$mpesa.Range("C11").Value = 173
$mpesa.Range("C12").Value = 92
$mpesa.Range("D12").Value = 910
$mpesa.Range("B20").Value = "KES 17,399"

# ---------------------------------------------------------------------
# Sheet: Unpaid Labor - insert the Fri 24/10 + Sat 25/10 rows before the
# "Total Unpaid Labor" row (old row 26 shifts down to row 35).
# ---------------------------------------------------------------------
$labor = $wb.Worksheets.Item("Unpaid Labor")
$labor.Range("A25:A33").EntireRow.Insert()

# Clone formatting from the existing data rows (4-24, style 5 / 10) onto
# the nine freshly inserted rows.
$labor.Range("A24:D24").Copy()
$labor.Range("A25:D33").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$laborRows = @(
    @(25, "24/10/2025", "Jack - UNPAID", 1500),
    @(26, "24/10/2025", "Fundi 1 - UNPAID", 1300),
    @(27, "24/10/2025", "Fundi 2 - UNPAID", 1300),
    @(28, "24/10/2025", "2 helpers @ 600 each - UNPAID", 1200),
    @(29, "25/10/2025", "Jack - UNPAID", 1500),
    @(30, "25/10/2025", "Fundi 1 - UNPAID", 1300),
    @(31, "25/10/2025", "Fundi 2 - UNPAID", 1300),
    @(32, "25/10/2025", "2 helpers @ 600 each - UNPAID", 1200),
    @(33, "25/10/2025", "Transport - UNPAID", 600)
)
foreach ($entry in $laborRows) {
    $r = $entry[0]
    $labor.Cells.Item($r, 1).Value = $entry[1]
    $labor.Cells.Item($r, 2).Value = $entry[2]
    $labor.Cells.Item($r, 3).Value = $entry[3]
    $labor.Cells.Item($r, 4).Value = "PENDING"
}

# Row 34 is left blank (matches the source layout); the old totals row
# (previously row 26) is now row 35.
$labor.Range("C35").Value = "KES 38,900"

# ---------------------------------------------------------------------
# Sheet: Pending Purchases
# ---------------------------------------------------------------------
$pending = $wb.Worksheets.Item("Pending Purchases")
$pending.Range("C14").Value = 37905
$pending.Range("C16").Value = "KES 170,405"
$pending.Range("B19").Value = "KES 1,516,311"
$pending.Range("B21").Value = "KES 38,900"
$pending.Range("B22").Value = "KES 170,405"
$pending.Range("B24").Value = "KES 1,806,916"
$pending.Range("B26").Value = "KES 806,916"
